# "Fixed qty of 1x2 headers"
# Row 8 is the "SIL VERTICAL PC TAIL PIN HEADER 1x2" (JP26) BOM line.
# Its Quantity (F8) was wrong; correct it to 6. The Total Cost formula in
# G8 (=F8*E8) recalculates automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value = 6

# Reflect the cell the editor left selected when saving.
$ws.Range("F8").Select() | Out-Null
